$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "87811004_0322_RO"
$ws.Range("B2").Value = 423
$ws.Range("D2").Value = 9983.42

$ws.Range("A3").Value = "87811004_0322_PE"
$ws.Range("B3").Value = 15
$ws.Range("D3").Value = 141.96

$ws.Range("A4").Value = "87811004_0322_HU"
$ws.Range("B4").Value = 954
$ws.Range("D4").Value = 1882603

$ws.Range("A5").Value = "87811004_0322_EU"
$ws.Range("B5").Value = 716
$ws.Range("D5").Value = 3090.3

$ws.Range("A6").Value = "87811004_0322_MX"
$ws.Range("B6").Value = 75
$ws.Range("D6").Value = 4825.1

$ws.Range("A7").Value = "87811004_0322_LL"
$ws.Range("B7").Value = 37
$ws.Range("D7").Value = 96.6

$ws.Range("A8").Value = "87811004_0322_BG"
$ws.Range("B8").Value = 6
$ws.Range("D8").Value = 22.44

$ws.Range("A9").Value = "87811004_0322_BR"
$ws.Range("B9").Value = 37
$ws.Range("D9").Value = 334.81

$ws.Range("A10").Value = "87811004_0322_CA"
$ws.Range("B10").Value = 324
$ws.Range("D10").Value = 1753.5

$ws.Range("A11").Value = "87811004_0322_CZ"
$ws.Range("B11").Value = 12
$ws.Range("D11").Value = 636.36

$ws.Range("A12").Value = "87811004_0322_CL"
$ws.Range("B12").Value = 31
$ws.Range("D12").Value = 63087

$ws.Range("A13").Value = "87811004_0322_CO"
$ws.Range("B13").Value = 31
$ws.Range("D13").Value = 286230

$ws.Range("A14").Value = "87811004_0322_NZ"
$ws.Range("B14").Value = 51
$ws.Range("D14").Value = 208.33

$ws.Range("A15").Value = "87811004_0322_AU"
$ws.Range("B15").Value = 367
$ws.Range("D15").Value = 3178.64

$ws.Range("A16").Value = "87811004_0322_CH"
$ws.Range("B16").Value = 79
$ws.Range("D16").Value = 243.39

$ws.Range("A17").Value = "87811004_0322_NO"
$ws.Range("B17").Value = 26
$ws.Range("D17").Value = 682.5

$ws.Range("A18").Value = "87811004_0322_US"
$ws.Range("B18").Value = 1487
$ws.Range("D18").Value = 10299.1

$ws.Range("A19").Value = "87811004_0322_DK"
$ws.Range("B19").Value = 21
$ws.Range("D19").Value = 268.24

$ws.Range("A20").Value = "87811004_0322_PL"
$ws.Range("B20").Value = 46
$ws.Range("D20").Value = 517.94

$ws.Range("A21").Value = "87811004_0322_SE"
$ws.Range("B21").Value = 32
$ws.Range("D21").Value = 931.77

$ws.Range("A22").Value = "87811004_0322_JP"
$ws.Range("B22").Value = 23
$ws.Range("D22").Value = 5950

$ws.Range("A23").Value = "87811004_0322_GB"
$ws.Range("B23").Value = 477
$ws.Range("D23").Value = 1593.33

